# Update Balance Sheet values on the "DHI" worksheet to reflect the
# restated figures for Inventory (row 4), Accounts Payable (row 15) and
# Long Term Tax Liability / Deferred (row 21) across columns B:F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DHI")

# Row 4 - Inventory
$ws.Range("B4").Value = 14476000000.0
$ws.Range("C4").Value = 13577000000.0
$ws.Range("D4").Value = 12237000000.0
$ws.Range("E4").Value = 12139000000.0
$ws.Range("F4").Value = 12225000000.0

# Row 15 - Accounts Payable
$ws.Range("B15").Value = 1137000000.0
$ws.Range("C15").Value = 846000000.0
$ws.Range("D15").Value = 901000000.0
$ws.Range("E15").Value = 761000000.0
$ws.Range("F15").Value = 707000000.0

# Row 21 - Long Term Tax Liability (Deferred)
$ws.Range("B21").Value = -143000000.0
$ws.Range("C21").Value = -142000000.0
$ws.Range("D21").Value = -145000000.0
$ws.Range("E21").Value = -159000000.0
$ws.Range("F21").Value = -150000000.0
